$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETS")

# Row 4 ("nuclear"): change all shareweights from 1 to 0 for the full
# year range. Update the year columns (C:AG) first so the dependent
# B4 formula recalculates to 0, then replace B4 itself with a literal
# 0 (it no longer holds a formula in the edited workbook).
$ws.Range("C4:AG4").Value = 0
$ws.Range("B4").Value = 0

# Row 17 ("municipal solid waste"): change all shareweights from 1 to 0.
# C17:AG17 become literal zeros; B17 keeps its "=C17" formula, which
# recalculates to 0 automatically.
$ws.Range("C17:AG17").Value = 0

# Reflect the last on-screen selection recorded in the saved workbook.
$ws.Range("A16").Select()

$wb.Save()
